$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Heading paragraph:
#    "RESPECTIVE CONTRIBUTIONS (one page limit)" -> "E. RESPECTIVE CONTRIBUTIONS"
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1Body = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$p1Body.Delete()

$p1 = $d.Paragraphs(1)
$insAt0 = $d.Range($p1.Range.Start, $p1.Range.Start)
$insAt0.InsertBefore('E. ')

$p1 = $d.Paragraphs(1)
$insAfterE = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$insAfterE.InsertAfter('RESPECTIVE CONTRIBUTIONS')

# ------------------------------------------------------------------
# 2) Body paragraph (3rd paragraph): replace the placeholder prompt
#    text with the actual respective-contributions narrative, and
#    drop the small (10pt) font formatting that wrapped it.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3Inner = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$p3Inner.Delete()

$p3 = $d.Paragraphs(3)
$insRun1 = $d.Range($p3.Range.Start, $p3.Range.Start)
$insRun1.InsertAfter('The proposed training plan was developed in collaboration with my research sponsors, Drs. Viola Vaccarino, Dr. Alvaro Alonso, and Dr. Amit J. Shah. In discussions with the sponsors, we determined the core training goals to be an expansion of the TL1 award and the MSCR, with a focus on advanced epidemiological methods, biostatistical/computational techniques, and foundational knowledge in autonomic control of cardiac physiology. We identified appropriate coursework offered at the Rollins School of Public Health, along with connections with collaborators (Drs. Jeanie Park and Marc Thames) to achieve these goals, along with experiential/mentored hands-on training through the research proposal. We developed the pilot study, with support from my advisors and collaborators, including Dr. Arshed Quyyumi, to study patients in the Emory Cardiovascular Biobank. ')

$p3 = $d.Paragraphs(3)
$afterRun1 = $p3.Range.End - 1
$insRun2 = $d.Range($afterRun1, $afterRun1)
$insRun2.InsertAfter('The participants have been continually enrolled by the research staff since the initiation of the pilot study in October 2019, and thus by the beginning of the award (as early as July 2019), we expect to have over 150 participants. I will perform the data analyses proposed in this training plan under, building upon the biostatistical training already received. I will work closely with Dr. Vaccarino to design the appropriate analyses to study the depression and heart rate variability, as she is an expert in psychological stress and its role in the development and prognosis of cardiovascular disease. I will work closely with Dr. Shah to extract and process ECG signal, as he is a cardiologist with expertise in autonomic function and the electrophysiology of stress. He will supervise the manual extraction of ECG data, adjudication of the signal quality, and the HRV analysis using the toolbox developed under his guidance. I will work closely with Dr. Alonso to model potential confounders and interactions. He has an expertise in cardiovascular epidemiology, and with his supervision I will expand upon my understanding of study design and epidemiological modeling. We will meet at least weekly with all the sponsors, although we will meet informally several times per week due to the proximity of our offices, the schedule lab meetings, and training lectures that are well attended by the entire EPICORE group. ')

$p3 = $d.Paragraphs(3)
$afterRun2 = $p3.Range.End - 1
$insRun3 = $d.Range($afterRun2, $afterRun2)
$insRun3.InsertAfter('Dr. Vaccarino will provide her expertise in psychological stress and research methods to help with interpretation of results. Dr. Alonso will provide his expertise in biostatistical modeling to help identify and understand potential confounders, mediators and interactions. Dr. Shah will provide his expertise in ECG signal processing to help troubleshoot technical aspects of the proposed analyses.')

# bookmark ("_GoBack") sits right after the 2nd run, before the 3rd.
# NB: placing a bookmark at a position that is the *very last* character
# offset of a paragraph (paragraph.End - 1) is unreliable, so the 3rd
# run's text is inserted first (above) and the bookmark is added
# afterwards at the (now interior) boundary between run 2 and run 3.
$bookmarkPos = $d.Range($afterRun2, $afterRun2)
$d.Bookmarks.Add('_GoBack', $bookmarkPos)

# ------------------------------------------------------------------
# 3) Add a trailing empty paragraph after the body paragraph.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphAfter()
